$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 843.0345
$ws.Cells.Item(129, 9).Value = 696.7778
$ws.Cells.Item(129, 10).Value = 908.85
$ws.Cells.Item(129, 11).Value = 2090.3334
$ws.Cells.Item(129, 12).Value = 2726.55
$ws.Cells.Item(129, 13).Value = 2909.6666
$ws.Cells.Item(129, 14).Value = -12726.55
$ws.Cells.Item(132, 8).Value = 3362.342
$ws.Cells.Item(132, 9).Value = 3965.9666
$ws.Cells.Item(132, 10).Value = 1098.75
$ws.Cells.Item(132, 11).Value = 11897.8998
$ws.Cells.Item(132, 12).Value = 3296.25
$ws.Cells.Item(132, 13).Value = -9367.899800000001
$ws.Cells.Item(132, 14).Value = -8356.25
$ws.Cells.Item(135, 8).Value = 11043057
$ws.Cells.Item(135, 9).Value = 580.7586
$ws.Cells.Item(135, 11).Value = 5226.8274
$ws.Cells.Item(135, 13).Value = -2691.8274
$ws.Cells.Item(137, 8).Value = 22728546
$ws.Cells.Item(137, 9).Value = 1227.8276
$ws.Cells.Item(137, 10).Value = 66668028
$ws.Cells.Item(137, 11).Value = 3683.4828
$ws.Cells.Item(137, 12).Value = 200004084
$ws.Cells.Item(137, 13).Value = -1133.4828
$ws.Cells.Item(137, 14).Value = -200009184
$ws.Cells.Item(138, 8).Value = 2836.6943
$ws.Cells.Item(138, 9).Value = 2471.6155
$ws.Cells.Item(138, 10).Value = 3043.0435
$ws.Cells.Item(138, 11).Value = 7414.8465
$ws.Cells.Item(138, 12).Value = 9129.130500000001
$ws.Cells.Item(138, 13).Value = -2274.8465
$ws.Cells.Item(138, 14).Value = -19409.1305

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 892535.5600000001
$ws.Cells.Item(2, 9).Value = 1127.619
$ws.Cells.Item(2, 10).Value = 2452499.5
$ws.Cells.Item(2, 11).Value = 1127.619
$ws.Cells.Item(2, 12).Value = 2452499.5
$ws.Cells.Item(2, 13).Value = -1014.619
$ws.Cells.Item(2, 14).Value = -2452725.5
$ws.Cells.Item(32, 8).Value = 1358.02
$ws.Cells.Item(32, 9).Value = 957.49414
$ws.Cells.Item(32, 10).Value = 3627.6667
$ws.Cells.Item(32, 11).Value = 957.49414
$ws.Cells.Item(32, 12).Value = 3627.6667
$ws.Cells.Item(32, 13).Value = -670.49414
$ws.Cells.Item(32, 14).Value = -4201.6667
$ws.Cells.Item(45, 8).Value = 1089.091
$ws.Cells.Item(45, 9).Value = 877.4
$ws.Cells.Item(45, 10).Value = 1151.3529
$ws.Cells.Item(45, 11).Value = 877.4
$ws.Cells.Item(45, 12).Value = 1151.3529
$ws.Cells.Item(45, 13).Value = -500.4
$ws.Cells.Item(45, 14).Value = -1905.3529
$ws.Cells.Item(61, 8).Value = 1950188.2
$ws.Cells.Item(61, 9).Value = 2179508
$ws.Cells.Item(61, 10).Value = 971
$ws.Cells.Item(61, 11).Value = 2179508
$ws.Cells.Item(61, 12).Value = 971
$ws.Cells.Item(61, 13).Value = -2179296
$ws.Cells.Item(61, 14).Value = -1395
$ws.Cells.Item(74, 8).Value = 8067435.5
$ws.Cells.Item(74, 9).Value = 9434498
$ws.Cells.Item(74, 10).Value = 16955.889
$ws.Cells.Item(74, 11).Value = 9434498
$ws.Cells.Item(74, 12).Value = 16955.889
$ws.Cells.Item(74, 13).Value = -9433624
$ws.Cells.Item(74, 14).Value = -18703.889
$ws.Cells.Item(77, 8).Value = 8067435.5
$ws.Cells.Item(77, 9).Value = 9434498
$ws.Cells.Item(77, 10).Value = 16955.889
$ws.Cells.Item(77, 11).Value = 47172490
$ws.Cells.Item(77, 12).Value = 84779.44499999999
$ws.Cells.Item(77, 13).Value = -47168122
$ws.Cells.Item(77, 14).Value = -93515.44499999999
$ws.Cells.Item(110, 8).Value = 1359.3182
$ws.Cells.Item(110, 9).Value = 1125.4117
$ws.Cells.Item(110, 11).Value = 1125.4117
$ws.Cells.Item(110, 13).Value = 919.5882999999999
$ws.Cells.Item(116, 8).Value = 892535.5600000001
$ws.Cells.Item(116, 9).Value = 1127.619
$ws.Cells.Item(116, 10).Value = 2452499.5
$ws.Cells.Item(116, 11).Value = 1127.619
$ws.Cells.Item(116, 12).Value = 2452499.5
$ws.Cells.Item(116, 13).Value = 1166.381
$ws.Cells.Item(116, 14).Value = -2457087.5
$ws.Cells.Item(122, 8).Value = 1360
$ws.Cells.Item(122, 9).Value = 1387.7894
$ws.Cells.Item(122, 10).Value = 1301.3334
$ws.Cells.Item(122, 11).Value = 4163.3682
$ws.Cells.Item(122, 12).Value = 3904.0002
$ws.Cells.Item(122, 13).Value = -1713.3682
$ws.Cells.Item(122, 14).Value = -8804.0002
$ws.Cells.Item(132, 8).Value = 7162077.5
$ws.Cells.Item(132, 9).Value = 8731044
$ws.Cells.Item(132, 10).Value = 101728.9
$ws.Cells.Item(132, 11).Value = 26193132
$ws.Cells.Item(132, 12).Value = 305186.7
$ws.Cells.Item(132, 13).Value = -26190602
$ws.Cells.Item(132, 14).Value = -310246.7
$ws.Cells.Item(136, 8).Value = 1950188.2
$ws.Cells.Item(136, 9).Value = 2179508
$ws.Cells.Item(136, 10).Value = 971
$ws.Cells.Item(136, 11).Value = 6538524
$ws.Cells.Item(136, 12).Value = 2913
$ws.Cells.Item(136, 13).Value = -6535974
$ws.Cells.Item(136, 14).Value = -8013

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 892535.5600000001
$ws.Cells.Item(3, 9).Value = 1127.619
$ws.Cells.Item(3, 10).Value = 2452499.5
$ws.Cells.Item(3, 11).Value = 1127.619
$ws.Cells.Item(3, 12).Value = 2452499.5
$ws.Cells.Item(3, 13).Value = -1013.619
$ws.Cells.Item(3, 14).Value = -2452727.5
$ws.Cells.Item(134, 8).Value = 3528910.2
$ws.Cells.Item(134, 9).Value = 4372524.5
$ws.Cells.Item(134, 10).Value = 1069.3636
$ws.Cells.Item(134, 11).Value = 13117573.5
$ws.Cells.Item(134, 12).Value = 3208.0908
$ws.Cells.Item(134, 13).Value = -13115038.5
$ws.Cells.Item(134, 14).Value = -8278.0908
$ws.Cells.Item(135, 8).Value = 32234
$ws.Cells.Item(135, 10).Value = 32234
$ws.Cells.Item(135, 12).Value = 32234
$ws.Cells.Item(135, 14).Value = -42374

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3928129.2
$ws.Cells.Item(31, 9).Value = 1115.6389
$ws.Cells.Item(31, 10).Value = 16780174
$ws.Cells.Item(31, 11).Value = 1115.6389
$ws.Cells.Item(31, 12).Value = 16780174
$ws.Cells.Item(31, 13).Value = -820.6388999999999
$ws.Cells.Item(31, 14).Value = -16780764
$ws.Cells.Item(34, 8).Value = 3928129.2
$ws.Cells.Item(34, 9).Value = 1115.6389
$ws.Cells.Item(34, 10).Value = 16780174
$ws.Cells.Item(34, 11).Value = 1115.6389
$ws.Cells.Item(34, 12).Value = 16780174
$ws.Cells.Item(34, 13).Value = -913.6388999999999
$ws.Cells.Item(34, 14).Value = -16780578
$ws.Cells.Item(58, 8).Value = 1030.3523
$ws.Cells.Item(58, 9).Value = 746.95776
$ws.Cells.Item(58, 10).Value = 2213.9412
$ws.Cells.Item(58, 11).Value = 746.95776
$ws.Cells.Item(58, 12).Value = 2213.9412
$ws.Cells.Item(58, 13).Value = -543.95776
$ws.Cells.Item(58, 14).Value = -2619.9412
$ws.Cells.Item(134, 8).Value = 1455.4117
$ws.Cells.Item(134, 9).Value = 1528.5518
$ws.Cells.Item(134, 11).Value = 4585.6554
$ws.Cells.Item(134, 13).Value = -2050.6554
$ws.Cells.Item(136, 8).Value = 1030.3523
$ws.Cells.Item(136, 9).Value = 746.95776
$ws.Cells.Item(136, 10).Value = 2213.9412
$ws.Cells.Item(136, 11).Value = 2240.87328
$ws.Cells.Item(136, 12).Value = 6641.823600000001
$ws.Cells.Item(136, 13).Value = 309.1267200000002
$ws.Cells.Item(136, 14).Value = -11741.8236

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 7463389
$ws.Cells.Item(5, 9).Value = 301.2069
$ws.Cells.Item(5, 11).Value = 903.6207000000001
$ws.Cells.Item(5, 13).Value = -791.6207000000001
$ws.Cells.Item(110, 8).Value = 2700
$ws.Cells.Item(110, 10).Value = 3500
$ws.Cells.Item(110, 12).Value = 10500
$ws.Cells.Item(110, 14).Value = -18680
$ws.Cells.Item(122, 8).Value = 10823792
$ws.Cells.Item(122, 9).Value = 20834008
$ws.Cells.Item(122, 11).Value = 187506072
$ws.Cells.Item(122, 13).Value = -187503622
$ws.Cells.Item(131, 8).Value = 6222.15
$ws.Cells.Item(131, 9).Value = 7901.25
$ws.Cells.Item(131, 10).Value = 5102.75
$ws.Cells.Item(131, 11).Value = 23703.75
$ws.Cells.Item(131, 12).Value = 15308.25
$ws.Cells.Item(131, 13).Value = -18663.75
$ws.Cells.Item(131, 14).Value = -25388.25
$ws.Cells.Item(132, 8).Value = 125001310
$ws.Cells.Item(132, 9).Value = 200000400
$ws.Cells.Item(132, 11).Value = 1800003600
$ws.Cells.Item(132, 13).Value = -1800001070
$ws.Cells.Item(135, 8).Value = 7463389
$ws.Cells.Item(135, 9).Value = 301.2069
$ws.Cells.Item(135, 11).Value = 2710.8621
$ws.Cells.Item(135, 13).Value = -175.8621000000003

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 21278524
$ws.Cells.Item(132, 9).Value = 30304800
$ws.Cells.Item(132, 10).Value = 2301.8572
$ws.Cells.Item(132, 11).Value = 90914400
$ws.Cells.Item(132, 12).Value = 6905.571599999999
$ws.Cells.Item(132, 13).Value = -90911870
$ws.Cells.Item(132, 14).Value = -11965.5716

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 2424.7188
$ws.Cells.Item(132, 9).Value = 2306.2952
$ws.Cells.Item(132, 10).Value = 4832.6665
$ws.Cells.Item(132, 11).Value = 6918.8856
$ws.Cells.Item(132, 12).Value = 14497.9995
$ws.Cells.Item(132, 13).Value = -4388.8856
$ws.Cells.Item(132, 14).Value = -19557.9995
$ws.Cells.Item(136, 8).Value = 901.7059
$ws.Cells.Item(136, 9).Value = 536.1875
$ws.Cells.Item(136, 10).Value = 6750
$ws.Cells.Item(136, 11).Value = 1608.5625
$ws.Cells.Item(136, 12).Value = 20250
$ws.Cells.Item(136, 13).Value = 941.4375
$ws.Cells.Item(136, 14).Value = -25350

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 931587.4
$ws.Cells.Item(14, 9).Value = 8207.143
$ws.Cells.Item(14, 10).Value = 2547502.8
$ws.Cells.Item(14, 11).Value = 8207.143
$ws.Cells.Item(14, 12).Value = 2547502.8
$ws.Cells.Item(14, 13).Value = -8039.143
$ws.Cells.Item(14, 14).Value = -2547838.8
$ws.Cells.Item(117, 8).Value = 22349.75
$ws.Cells.Item(117, 10).Value = 22349.75
$ws.Cells.Item(117, 12).Value = 22349.75
$ws.Cells.Item(117, 14).Value = -31527.75
$ws.Cells.Item(132, 8).Value = 6764613.5
$ws.Cells.Item(132, 9).Value = 7133408.5
$ws.Cells.Item(132, 10).Value = 3367.6667
$ws.Cells.Item(132, 11).Value = 21400225.5
$ws.Cells.Item(132, 12).Value = 10103.0001
$ws.Cells.Item(132, 13).Value = -21397695.5
$ws.Cells.Item(132, 14).Value = -15163.0001
$ws.Cells.Item(136, 8).Value = 2806119
$ws.Cells.Item(136, 9).Value = 7298.4194
$ws.Cells.Item(136, 10).Value = 7144291
$ws.Cells.Item(136, 11).Value = 21895.2582
$ws.Cells.Item(136, 12).Value = 21432873
$ws.Cells.Item(136, 13).Value = -19345.2582
$ws.Cells.Item(136, 14).Value = -21437973
